# Update scripts with new tpm:
# The underlying TPM values were recomputed and the "ECs" sending-cluster
# rows (old rows 2-4) were dropped, leaving only the "FAPs" and "MuSCs"
# sending-cluster blocks (which shift up to rows 2-4 and 5-7), with their
# numeric columns refreshed to reflect the new TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three "ECs" sending-cluster rows that used to sit at rows 2-4;
# deleting the same row index three times in a row removes them all and
# shifts everything below up, shrinking the sheet to 7 rows (A1:T7).
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Refreshed data (sending cluster, ligand, receptor, target cluster, then
# the recalculated numeric columns E through T) for the remaining rows.
$newRows = @(
    @('FAPs',  'Efna2', 'Epha1', 'ECs',   3, 1, 2.015377,             6.046131,             0.7554960962715589, 0.7554960962715588,
      3, 1, 2.211928,             6.635783999999999,  0.174938892641363,   0.1749388926413629,
      4.457868816856,        40.12081935170399,    0.1321656504766191,  0.132165650476619),
    @('FAPs',  'Efna2', 'Epha1', 'FAPs',  3, 1, 2.015377,             6.046131,             0.7554960962715589, 0.7554960962715588,
      3, 1, 3.864911333333334,   11.594734,           0.3056714815357404,  0.3056714815357404,
      7.789253408239333,     70.103280674154,      0.2309336110417957,  0.2309336110417957),
    @('FAPs',  'Efna2', 'Epha1', 'MuSCs', 3, 1, 2.015377,             6.046131,             0.7554960962715589, 0.7554960962715588,
      3, 1, 6.567164333333333,   19.701493,           0.5193896258228967,  0.5193896258228966,
      13.23531195262033,     119.117807573583,     0.3923968347531441,  0.392396834753144),
    @('MuSCs', 'Efna2', 'Epha1', 'ECs',   3, 1, 0.6522436666666667,   1.956731,             0.2445039037284412, 0.2445039037284411,
      3, 1, 2.211928,             6.635783999999999,  0.174938892641363,   0.1749388926413629,
      1.442716029122667,     12.984444262104,      0.04277324216474392, 0.0427732421647439),
    @('MuSCs', 'Efna2', 'Epha1', 'FAPs',  3, 1, 0.6522436666666667,   1.956731,             0.2445039037284412, 0.2445039037284411,
      3, 1, 3.864911333333334,   11.594734,           0.3056714815357404,  0.3056714815357404,
      2.520863939394889,     22.687775454554,      0.07473787049394465, 0.07473787049394463),
    @('MuSCs', 'Efna2', 'Epha1', 'MuSCs', 3, 1, 0.6522436666666667,   1.956731,             0.2445039037284412, 0.2445039037284411,
      3, 1, 6.567164333333333,   19.701493,           0.5193896258228967,  0.5193896258228966,
      4.283391344375889,     38.550522099383,      0.1269927910697526,  0.1269927910697526)
)

$rowIndex = 2
foreach ($rowData in $newRows) {
    $colIndex = 1
    foreach ($value in $rowData) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
